$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New entry for the journal: date, description of today's activity, duration.
$newRow = 46
$newDate = 43174
$newText = "J'ai créé un dossier qui contiendra mon application pour le téléphone sur le site de Ionic directement. Ce qui me permettra de faire mon développement sur l'ordinateur et de voir sur mon natel mon application. J'ai regardé avec M. Carrel avec quoi je devais démarrer et je dois commencer à implémenter la lecture des QR Code."
$newDuration = "4 périodes"

# Reuse the formatting of the row above (date style, wrap-text style, plain style)
# so no new cell styles are minted in styles.xml.
$ws.Range("A45").Copy()
$ws.Range("A46").PasteSpecial(-4122)

$ws.Range("B45").Copy()
$ws.Range("B46").PasteSpecial(-4122)

$ws.Range("C45").Copy()
$ws.Range("C46").PasteSpecial(-4122)

$ws.Range("A46").Value = $newDate
$ws.Range("B46").Value = $newText
$ws.Range("C46").Value = $newDuration

# Match the row height Excel computed for the wrapped text of the new entry.
$ws.Rows.Item($newRow).RowHeight = 60

# Move the active selection one row past the newly added entry, same as Excel
# does after the user finishes typing in the last cell of the new row.
$ws.Range("C47").Select()
